# bioSample_2002.xlsx — "updated files with strain names"
#
# The harvester column (B) for every data row (2-49) was changed from the
# shared string "S.GISH" to a new shared string "H.BROWN". Writing the new
# literal value through Range.Value lets the engine's own shared-string
# table management retire the now-unused "S.GISH" entry and append
# "H.BROWN" at the end, which is exactly the sharedStrings.xml shape the
# diff shows (every index >= the old "S.GISH" slot shifts down by one, and
# "H.BROWN" lands as the new last entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 2).Value = "H.BROWN"
}

# Selection moved from F3:F4 (active F3) to B3:B49 (active B3).
[void]$ws.Range("B3:B49").Select()
